# Adds a "Conclusion" section to the end of the report.
#
# The document's final paragraph holds only the (hidden) "_GoBack"
# bookmark left over from the last save. We replace that paragraph's
# range with a heading paragraph ("Conclusion", styled like the other
# section headings in the doc) followed by the conclusion body
# paragraph, keeping the _GoBack bookmark mid-sentence exactly where
# Word itself would have left it after the author's final keystrokes.

$d = $word.ActiveDocument

$conclusionXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>Conclusion</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:t xml:space="preserve">In conclusion, </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>GitHub</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> is an essential platform for modern software development, offering a complete environment for version control, collaboration, and automation. Its combination of </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>Git</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>-</w:t>
  </w:r>
  <w:r>
    <w:lastRenderedPageBreak/>
    <w:t xml:space="preserve">based workflow, strong community support, and integrated CI/CD features makes it one of the most reliable and widely adopted tools for developers. By using </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>GitHub</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>, teams can improve productivity, maintain clean ve</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r>
    <w:t>rsion histories, and deliver high-quality software efficiently.</w:t>
  </w:r>
</w:p>
"@

$lastParagraph = $d.Paragraphs.Last
$lastParagraph.Range.InsertXML($conclusionXml) | Out-Null
